# Auto-generated edit script: updates cached market-price / profit values
# in several Leve-profit worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 398
$ws.Range("I11").Value = 398
$ws.Range("K11").Value = 398
$ws.Range("M11").Value = -258
$ws.Range("H29").Value = 1300
$ws.Range("I29").Value = 1000
$ws.Range("J29").Value = 1800
$ws.Range("K29").Value = 3000
$ws.Range("L29").Value = 5400
$ws.Range("M29").Value = -2719
$ws.Range("N29").Value = -5962
$ws.Range("H70").Value = 3111
$ws.Range("I70").Value = 1566.6666
$ws.Range("J70").Value = 3772.8572
$ws.Range("K70").Value = 4699.9998
$ws.Range("L70").Value = 11318.5716
$ws.Range("M70").Value = -4429.9998
$ws.Range("N70").Value = -11858.5716
$ws.Range("H73").Value = 3111
$ws.Range("I73").Value = 1566.6666
$ws.Range("J73").Value = 3772.8572
$ws.Range("K73").Value = 4699.9998
$ws.Range("L73").Value = 11318.5716
$ws.Range("M73").Value = -3763.9998
$ws.Range("N73").Value = -13190.5716
$ws.Range("H80").Value = 947.05884
$ws.Range("J80").Value = 1466.5555
$ws.Range("L80").Value = 4399.666499999999
$ws.Range("N80").Value = -6395.666499999999
$ws.Range("H82").Value = 8013.8
$ws.Range("I82").Value = 947
$ws.Range("J82").Value = 11042.429
$ws.Range("K82").Value = 2841
$ws.Range("L82").Value = 33127.287
$ws.Range("M82").Value = -2435
$ws.Range("N82").Value = -33939.287
$ws.Range("H83").Value = 947.05884
$ws.Range("J83").Value = 1466.5555
$ws.Range("L83").Value = 13198.9995
$ws.Range("N83").Value = -23182.9995
$ws.Range("H85").Value = 8013.8
$ws.Range("I85").Value = 947
$ws.Range("J85").Value = 11042.429
$ws.Range("K85").Value = 2841
$ws.Range("L85").Value = 33127.287
$ws.Range("M85").Value = -1437
$ws.Range("N85").Value = -35935.287
$ws.Range("H88").Value = 3400
$ws.Range("I88").Value = 5000
$ws.Range("J88").Value = 2200
$ws.Range("K88").Value = 5000
$ws.Range("L88").Value = 2200
$ws.Range("M88").Value = -4594
$ws.Range("N88").Value = -3012
$ws.Range("H91").Value = 3400
$ws.Range("I91").Value = 5000
$ws.Range("J91").Value = 2200
$ws.Range("K91").Value = 5000
$ws.Range("L91").Value = 2200
$ws.Range("M91").Value = -3596
$ws.Range("N91").Value = -5008
$ws.Range("H109").Value = 34657.895
$ws.Range("J109").Value = 34657.895
$ws.Range("L109").Value = 34657.895
$ws.Range("N109").Value = -37431.895
$ws.Range("H137").Value = 1907250.5
$ws.Range("I137").Value = 2802756.5
$ws.Range("J137").Value = 4300.375
$ws.Range("K137").Value = 8408269.5
$ws.Range("L137").Value = 12901.125
$ws.Range("M137").Value = -8405719.5
$ws.Range("N137").Value = -18001.125
$ws.Range("H138").Value = 2449.4106
$ws.Range("I138").Value = 729.9583
$ws.Range("J138").Value = 3030.6338
$ws.Range("K138").Value = 2189.8749
$ws.Range("L138").Value = 9091.901400000001
$ws.Range("M138").Value = 2950.1251
$ws.Range("N138").Value = -19371.9014

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1376.86
$ws.Range("I32").Value = 1151.6282
$ws.Range("J32").Value = 2175.4092
$ws.Range("K32").Value = 1151.6282
$ws.Range("L32").Value = 2175.4092
$ws.Range("M32").Value = -864.6282000000001
$ws.Range("N32").Value = -2749.4092
$ws.Range("H132").Value = 2482.7368
$ws.Range("I132").Value = 1219.5714
$ws.Range("J132").Value = 6019.6
$ws.Range("K132").Value = 3658.7142
$ws.Range("L132").Value = 18058.8
$ws.Range("M132").Value = -1128.7142
$ws.Range("N132").Value = -23118.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3596.389
$ws.Range("I134").Value = 1670.8462
$ws.Range("J134").Value = 8602.799999999999
$ws.Range("K134").Value = 5012.5386
$ws.Range("L134").Value = 25808.4
$ws.Range("M134").Value = -2477.5386
$ws.Range("N134").Value = -30878.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6238.591
$ws.Range("I31").Value = 2141.72
$ws.Range("K31").Value = 2141.72
$ws.Range("M31").Value = -1846.72
$ws.Range("H34").Value = 6238.591
$ws.Range("I34").Value = 2141.72
$ws.Range("K34").Value = 2141.72
$ws.Range("M34").Value = -1939.72
$ws.Range("H62").Value = 2901.6667
$ws.Range("I62").Value = 2952.5
$ws.Range("J62").Value = 2800
$ws.Range("K62").Value = 2952.5
$ws.Range("L62").Value = 2800
$ws.Range("M62").Value = -2328.5
$ws.Range("N62").Value = -4048
$ws.Range("H65").Value = 2901.6667
$ws.Range("I65").Value = 2952.5
$ws.Range("J65").Value = 2800
$ws.Range("K65").Value = 14762.5
$ws.Range("L65").Value = 14000
$ws.Range("M65").Value = -11642.5
$ws.Range("N65").Value = -20240
$ws.Range("H99").Value = 13337627
$ws.Range("I99").Value = 33335000
$ws.Range("K99").Value = 33335000
$ws.Range("M99").Value = -33333502
$ws.Range("H126").Value = 13337627
$ws.Range("I126").Value = 33335000
$ws.Range("K126").Value = 100005000
$ws.Range("M126").Value = -100002530

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1028602.56
$ws.Range("I5").Value = 616.6667
$ws.Range("J5").Value = 1909733.2
$ws.Range("K5").Value = 1850.0001
$ws.Range("L5").Value = 5729199.6
$ws.Range("M5").Value = -1738.0001
$ws.Range("N5").Value = -5729423.6
$ws.Range("H122").Value = 2667.9443
$ws.Range("I122").Value = 1063.5555
$ws.Range("J122").Value = 2988.8223
$ws.Range("K122").Value = 9571.9995
$ws.Range("L122").Value = 26899.4007
$ws.Range("M122").Value = -7121.9995
$ws.Range("N122").Value = -31799.4007
$ws.Range("H135").Value = 1028602.56
$ws.Range("I135").Value = 616.6667
$ws.Range("J135").Value = 1909733.2
$ws.Range("K135").Value = 5550.0003
$ws.Range("L135").Value = 17187598.8
$ws.Range("M135").Value = -3015.0003
$ws.Range("N135").Value = -17192668.8
$ws.Range("H140").Value = 6021.6665
$ws.Range("I140").Value = 10743.333
$ws.Range("J140").Value = 1300
$ws.Range("K140").Value = 32229.999
$ws.Range("L140").Value = 3900
$ws.Range("M140").Value = -27049.999
$ws.Range("N140").Value = -14260

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4218.4116
$ws.Range("I122").Value = 3476.0833
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 10428.2499
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -7978.249899999999
$ws.Range("N122").Value = -22900
$ws.Range("H132").Value = 3582.5881
$ws.Range("I132").Value = 2148.6667
$ws.Range("J132").Value = 4364.727
$ws.Range("K132").Value = 6446.000100000001
$ws.Range("L132").Value = 13094.181
$ws.Range("M132").Value = -3916.000100000001
$ws.Range("N132").Value = -18154.181
$ws.Range("H141").Value = 31971.5
$ws.Range("J141").Value = 31165.8
$ws.Range("L141").Value = 31165.8
$ws.Range("N141").Value = -41525.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 49999.668
$ws.Range("J24").Value = 49999.668
$ws.Range("L24").Value = 49999.668
$ws.Range("N24").Value = -50685.668
$ws.Range("H25").Value = 15000
$ws.Range("J25").Value = 15000
$ws.Range("L25").Value = 15000
$ws.Range("N25").Value = -15460
$ws.Range("H122").Value = 5506.615
$ws.Range("I122").Value = 3842.3333
$ws.Range("J122").Value = 9251.25
$ws.Range("K122").Value = 11526.9999
$ws.Range("L122").Value = 27753.75
$ws.Range("M122").Value = -9076.999899999999
$ws.Range("N122").Value = -32653.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 17500
$ws.Range("J12").Value = 17500
$ws.Range("L12").Value = 17500
$ws.Range("N12").Value = -17784
$ws.Range("H126").Value = 2249.348
$ws.Range("I126").Value = 1076.3889
$ws.Range("K126").Value = 3229.1667
$ws.Range("M126").Value = -759.1666999999998

